$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K (strikeout) values for rows 2-21, replacing the old Strike# values
$kValues = @{
    2  = 3
    3  = 6
    4  = 0
    5  = 4
    6  = 7
    7  = 9
    8  = 5
    9  = 5
    10 = 5
    11 = 6
    12 = 6
    13 = 7
    14 = 6
    15 = 1
    16 = 2
    17 = 2
    18 = 7
    19 = 3
    20 = 3
    21 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
